$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 390 (this shifts existing rows 390+ down by 3,
# carrying along their formatting/values, just as Excel's Insert does).
$ws.Rows("390:392").Insert()

# The common ("header-like") columns are constant for every row in this block,
# so copy them from the row directly below (old row 390, now shifted to 393).
$commonCols = @("A","B","C","E","F","G","H","I","J","K","Q","R","T")
foreach ($col in $commonCols) {
    $srcValue = $ws.Range("$col`393").Value2
    $ws.Range("$col`390").Value = $srcValue
    $ws.Range("$col`391").Value = $srcValue
    $ws.Range("$col`392").Value = $srcValue
}

# New week: 2021-11-05 (serial 44505)
$ws.Range("D390").Value = 44505
$ws.Range("D391").Value = 44505
$ws.Range("D392").Value = 44505

# Row 390: Especial
$ws.Range("L390").Value = "Especial"
$ws.Range("M390").Value = 500
$ws.Range("N390").Value = 12000
$ws.Range("O390").Value = 12500
$ws.Range("P390").Value = 12250
$ws.Range("S390").Value = 1750

# Row 391: Primera
$ws.Range("L391").Value = "Primera"
$ws.Range("M391").Value = 600
$ws.Range("N391").Value = 10000
$ws.Range("O391").Value = 10500
$ws.Range("P391").Value = 10250
$ws.Range("S391").Value = 1464

# Row 392: Segunda
$ws.Range("L392").Value = "Segunda"
$ws.Range("M392").Value = 360
$ws.Range("N392").Value = 8000
$ws.Range("O392").Value = 8500
$ws.Range("P392").Value = 8250
$ws.Range("S392").Value = 1179
